# Update C3 with new CP2102 according to spec
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook-level view settings (best effort) ---
$excel.ActiveWindow.Width = 14400
$excel.ActiveWindow.Height = 6020

# --- A3: add C19 to the designator list ---
$ws.Range("A3").Value = "C2,C7,C9,C10,C19,C1,C3"

# --- Insert a new row before row 8 for R13, shifting old rows 8..28 down by one ---
$ws.Rows("8:8").Insert()
$ws.Rows(8).RowHeight = 13.5

# Row 8 (new): R13
$ws.Range("A8").Value = "R13"
$ws.Range("B8").Value = "R_0603_1608Metric"
$ws.Range("C8").Value = "1K"
$ws.Range("D8").Value = "C21190"

# --- Fix C13 comment text (now row 20 after insertion): "100 nF" -> "100nF" ---
$ws.Range("C20").Value = "100nF"

# --- U2 row (now row 26 after insertion): update footprint & comment for new CP2102 part ---
$ws.Range("B26").Value = "CP2102-QFN50P500X500X80-29N-D"
$ws.Range("C26").Value = "CP2102N-A02-GQFN28R"

# --- Insert a new row before row 27 (old row27 'C12' etc shift down) for R15 ---
$ws.Rows("27:27").Insert()
$ws.Rows(27).RowHeight = 13.5

# Row 27 (new): R15
$ws.Range("A27").Value = "R15"
$ws.Range("B27").Value = "R_0402_1005Metric"
$ws.Range("C27").Value = "47K"
$ws.Range("D27").Value = "C137974"

# --- Append new row 31 for R14 (after L1 which is now at row 30) ---
$ws.Rows(31).RowHeight = 13.5
$ws.Range("A31").Value = "R14"
$ws.Range("B31").Value = "R_0201_0603Metric"
$ws.Range("C31").Value = "22K"
$ws.Range("D31").Value = "C31850"

# --- Update sheet view: scroll position and selection ---
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D31").Select()
